$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Förändrad" (Changed) date column C, rows 2-43, advances by one day
# from 45755 (2025-04-08) to 45756 (2025-04-09).
$range = $ws.Range("C2:C43")
$range.Value = 45756
